# Regenerate merged AHB files
# - Rename the "_old"/"_new" header suffixes to the concrete version tags
#   (FV2410 / FV2504) used by the regenerated merge.
# - Turn the data range into a real Excel Table ("Table1").
# - Freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) ------------------------------------------------
# Columns A-J were "<Label>_old" -> "<Label>_FV2410"
# Column  K  is "diff" and is left untouched
# Columns L-U were "<Label>_new" -> "<Label>_FV2504"
$headerRenames = [ordered]@{
    "A1" = "Segmentname_FV2410"
    "B1" = "Segmentgruppe_FV2410"
    "C1" = "Segment_FV2410"
    "D1" = "Datenelement_FV2410"
    "E1" = "Segment ID_FV2410"
    "F1" = "Code_FV2410"
    "G1" = "Qualifier_FV2410"
    "H1" = "Beschreibung_FV2410"
    "I1" = "Bedingungsausdruck_FV2410"
    "J1" = "Bedingung_FV2410"
    "L1" = "Segmentname_FV2504"
    "M1" = "Segmentgruppe_FV2504"
    "N1" = "Segment_FV2504"
    "O1" = "Datenelement_FV2504"
    "P1" = "Segment ID_FV2504"
    "Q1" = "Code_FV2504"
    "R1" = "Qualifier_FV2504"
    "S1" = "Beschreibung_FV2504"
    "T1" = "Bedingungsausdruck_FV2504"
    "U1" = "Bedingung_FV2504"
}

foreach ($addr in $headerRenames.Keys) {
    $ws.Range($addr).Value = $headerRenames[$addr]
}

# --- 2. Turn A1:U54 into a native Excel table -------------------------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U54"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row ------------------------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "AHB workbook regenerated: header labels updated, Table1 created, header row frozen."
